$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# --- Row 1 : base operands -------------------------------------------------
$ws2.Range("A1").Value = 123
$ws2.Range("B1").Value = 456

# --- Column A labels (string order: Summe, Multiplikat, Division, Subtraktion)
$ws2.Range("A2").Value = "Summe"
$ws2.Range("B2").Formula = "=A1+B1"

$ws2.Range("A4").Value = "Multiplikat"
$ws2.Range("B4").Formula = "=A1*B1"

$ws2.Range("A5").Value = "Division"
$ws2.Range("B5").Formula = "=B1/A1"

$ws2.Range("A3").Value = "Subtraktion"
$ws2.Range("B3").Formula = "=A1-B1"

# --- Column C formatted results (numFmt creation order matches the fixture)
$ws2.Range("C3").NumberFormat = '#,##0.0\ "€"'
$ws2.Range("C3").Formula = "=B3"

$ws2.Range("C4").NumberFormat = '#,##0.000'
$ws2.Range("C4").Formula = "=B4"

$ws2.Range("C5").NumberFormat = '[$¥-411]#,##0.000'
$ws2.Range("C5").Formula = "=B5"

$ws2.Range("C6").NumberFormat = '#,##0.0\ [$₽-419]'
$ws2.Range("C6").Formula = "=B5"

$ws2.Range("C7").NumberFormat = "0.00%"
$ws2.Range("C7").Formula = "=B5"

$ws2.Range("C8").NumberFormat = '[$-F800]dddd\,\ mmmm\ dd\,\ yyyy'
$ws2.Range("C8").Value = 43544

$ws2.Range("C2").NumberFormat = "0.0E+00"
$ws2.Range("C2").Formula = "=B2"

# --- Row 1 header / merged "formatted" cell --------------------------------
$ws2.Range("C1:D1").Merge()
$ws2.Range("C1:D1").HorizontalAlignment = -4108
$ws2.Range("C1").Value = "formatted"

# --- Column D descriptions --------------------------------------------------
$ws2.Range("D2").Value = "Wissenschaftlich"
$ws2.Range("D4").Value = "Zahl - Tausenderzeichen - 3 Stellen"
$ws2.Range("D5").Value = "Währung - japanisch - 3 Stellen"
$ws2.Range("D6").Value = "Währung - russisch - 1 Stelle"
$ws2.Range("D7").Value = "Prozent - 2 Stellen"
$ws2.Range("D3").Value = "Währung - Euro - 1 Stelle"
$ws2.Range("D8").Value = "Datum"

# --- Row 9 : second scientific format --------------------------------------
$ws2.Range("C9").NumberFormat = "0.000E+00"
$ws2.Range("C9").Formula = "=B2"
$ws2.Range("D9").Value = "Wissenschaftlich #2"

# --- Column widths -----------------------------------------------------------
$ws2.Columns.Item(3).ColumnWidth = 22.42578125
$ws2.Columns.Item(4).ColumnWidth = 28.140625

# --- Page setup (paper size / orientation like sheet1) -----------------------
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

Write-Host "edit applied"
